# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker detail table (rows 16-31, columns C:G) is rebuilt: rows are
# regrouped per worker (sorted by descending "Periodo Mora"), and
# ELSA MARIA BARRIOS COVA's "Salario Basico" is corrected from 781242 to
# 689455 for all of her periods. Column B ("CC") and the header rows are
# unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1049452296"
$ws.Range("D16").Value = "YAMIRIS ESTER RODRIGUEZ PUELLO"
$ws.Range("E16").Value = "1808"
$ws.Range("F16").Value = 43916
$ws.Range("G16").Value = 1097092

$ws.Range("C17").Value = "1049452296"
$ws.Range("D17").Value = "YAMIRIS ESTER RODRIGUEZ PUELLO"
$ws.Range("E17").Value = "1807"
$ws.Range("F17").Value = 43916
$ws.Range("G17").Value = 1097092

$ws.Range("C18").Value = "44151246"
$ws.Range("D18").Value = "ELSA MARIA BARRIOS COVA"
$ws.Range("E18").Value = "1808"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 689455

$ws.Range("C19").Value = "44151246"
$ws.Range("D19").Value = "ELSA MARIA BARRIOS COVA"
$ws.Range("E19").Value = "1807"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 689455

$ws.Range("C20").Value = "44151246"
$ws.Range("D20").Value = "ELSA MARIA BARRIOS COVA"
$ws.Range("E20").Value = "1806"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 689455

$ws.Range("C21").Value = "44151246"
$ws.Range("D21").Value = "ELSA MARIA BARRIOS COVA"
$ws.Range("E21").Value = "1805"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 689455

$ws.Range("C22").Value = "44151246"
$ws.Range("D22").Value = "ELSA MARIA BARRIOS COVA"
$ws.Range("E22").Value = "1804"
$ws.Range("F22").Value = 31249
$ws.Range("G22").Value = 689455

$ws.Range("C23").Value = "44151246"
$ws.Range("D23").Value = "ELSA MARIA BARRIOS COVA"
$ws.Range("E23").Value = "1803"
$ws.Range("F23").Value = 31249
$ws.Range("G23").Value = 689455

$ws.Range("C24").Value = "73432535"
$ws.Range("D24").Value = "SAULO ELVIRO ROMERO CUETO"
$ws.Range("E24").Value = "1804"
$ws.Range("F24").Value = 10616
$ws.Range("G24").Value = 884667

$ws.Range("C25").Value = "73432535"
$ws.Range("D25").Value = "SAULO ELVIRO ROMERO CUETO"
$ws.Range("E25").Value = "1803"
$ws.Range("F25").Value = 35387
$ws.Range("G25").Value = 884667

$ws.Range("C26").Value = "73432535"
$ws.Range("D26").Value = "SAULO ELVIRO ROMERO CUETO"
$ws.Range("E26").Value = "1802"
$ws.Range("F26").Value = 35387
$ws.Range("G26").Value = 884667

$ws.Range("C27").Value = "1052078967"
$ws.Range("D27").Value = "ELVIA PATRICIA MERCADO ORTEGA"
$ws.Range("E27").Value = "1808"
$ws.Range("F27").Value = 31249
$ws.Range("G27").Value = 781242

$ws.Range("C28").Value = "1052078967"
$ws.Range("D28").Value = "ELVIA PATRICIA MERCADO ORTEGA"
$ws.Range("E28").Value = "1807"
$ws.Range("F28").Value = 31249
$ws.Range("G28").Value = 781242

$ws.Range("C29").Value = "1052094888"
$ws.Range("D29").Value = "NATHALY PAOLA DEL VALLE RODRIGUEZ"
$ws.Range("E29").Value = "1804"
$ws.Range("F29").Value = 9375
$ws.Range("G29").Value = 781242

$ws.Range("C30").Value = "1052094888"
$ws.Range("D30").Value = "NATHALY PAOLA DEL VALLE RODRIGUEZ"
$ws.Range("E30").Value = "1803"
$ws.Range("F30").Value = 31249
$ws.Range("G30").Value = 781242

$ws.Range("C31").Value = "1052094888"
$ws.Range("D31").Value = "NATHALY PAOLA DEL VALLE RODRIGUEZ"
$ws.Range("E31").Value = "1802"
$ws.Range("F31").Value = 31249
$ws.Range("G31").Value = 781242
